$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0911
$ws.Range("E2").Value = -0.00462
$ws.Range("F2").Value = 0.301
$ws.Range("G2").Value = 0.1287502195678904
$ws.Range("H2").Value = 0.1287502195678904
$ws.Range("I2").Value = 0.0420918383850227
$ws.Range("J2").Value = 0.03679994540887101
$ws.Range("K2").Value = 492.4
$ws.Range("L2").Value = 0.02162304584577551
$ws.Range("M2").Value = 489.2
$ws.Range("N2").Value = 0.02962352926928224
$ws.Range("O2").Value = 0.9935012185215273
$ws.Range("P2").Value = 178
$ws.Range("Q2").Value = 0.01077879846674619
$ws.Range("R2").Value = 0.3614947197400488
$ws.Range("S2").Value = 311.2
$ws.Range("T2").Value = 0.6361406377759607
$ws.Range("U2").Value = 4486.6
$ws.Range("V2").Value = 0.2716862764095701
$ws.Range("W2").Value = 0.0210256031470746
$ws.Range("X2").Value = 0.08526244966425546
$ws.Range("Y2").Value = -0.06423684651718087
$ws.Range("Z2").Value = 1.017237317675182
$ws.Range("AA2").Value = 0.02946221068873586
$ws.Range("AB2").Value = 0.0661207122965583
$ws.Range("AC2").Value = -0.03665850160782244
$ws.Range("AD2").Value = 6022.6
$ws.Range("AE2").Value = 338.4232814813149
$ws.Range("AF2").Value = 6361.023281481315
$ws.Range("AG2").Value = 1874.423281481315
$ws.Range("AH2").Value = 0.2780784531256095
$ws.Range("AI2").Value = 0.2239140697764415
$ws.Range("AJ2").Value = 0.1019355192307842
$ws.Range("AK2").Value = 0.07835653223747367
$ws.Range("AL2").Value = 259.5
$ws.Range("AM2").Value = 259.5
$ws.Range("AN2").Value = 4.232921000843408
$ws.Range("AO2").Value = 3.746435452793834
$ws.Range("AP2").Value = 1.317418668457489
$ws.Range("AQ2").Value = 3.746435452793834

# Row 3
$ws.Range("D3").Value = 0.0672
$ws.Range("E3").Value = -0.00462
$ws.Range("G3").Value = 0.1110619469026549
$ws.Range("H3").Value = 0.1110619469026549
$ws.Range("I3").Value = 0.06210048755781198
$ws.Range("J3").Value = 0.04648562507940261
$ws.Range("K3").Value = 518
$ws.Range("L3").Value = 0.03638151425762045
$ws.Range("M3").Value = 343
$ws.Range("N3").Value = 0.04356164035611323
$ws.Range("O3").Value = 0.6621621621621622
$ws.Range("P3").Value = 178
$ws.Range("Q3").Value = 0.022606332313085
$ws.Range("R3").Value = 0.3436293436293436
$ws.Range("S3").Value = 165
$ws.Range("T3").Value = 0.4810495626822158
$ws.Range("U3").Value = 3256
$ws.Range("V3").Value = 0.4135180787157571
$ws.Range("W3").Value = 0.04495084044187197
$ws.Range("X3").Value = 0.09056202210470354
$ws.Range("Y3").Value = -0.04561118166283157
$ws.Range("Z3").Value = 1.097775388027257
$ws.Range("AA3").Value = 0.05103077510923078
$ws.Range("AB3").Value = 0.0662189207618023
$ws.Range("AC3").Value = -0.01518814565257152
$ws.Range("AD3").Value = 3981
$ws.Range("AE3").Value = 64.06629075936495
$ws.Range("AF3").Value = 4045.066290759365
$ws.Range("AG3").Value = 789.0662907593651
$ws.Range("AH3").Value = 0.3393806301722204
$ws.Range("AI3").Value = 0.2343734141009195
$ws.Range("AJ3").Value = 0.0910850007117133
$ws.Range("AK3").Value = 0.05634953619266021
$ws.Range("AL3").Value = 170
$ws.Range("AM3").Value = 170
$ws.Range("AN3").Value = 3.478982784234903
$ws.Range("AO3").Value = 5.188235294117647
$ws.Range("AP3").Value = 0.6895624318442412
$ws.Range("AQ3").Value = 5.188235294117647

# Row 4
$ws.Range("D4").Value = 0.115
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 0.301
$ws.Range("G4").Value = 0.1582610733536442
$ws.Range("H4").Value = 0.1582610733536442
$ws.Range("I4").Value = 0.008709702584439889
$ws.Range("J4").Value = 0.008709702584439889
$ws.Range("K4").Value = -25.6
$ws.Range("L4").Value = -0.002999765643309117
$ws.Range("M4").Value = 146.2
$ws.Range("N4").Value = 0.0169212962962963
$ws.Range("O4").Value = -5.710937499999999
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 146.2
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 1230.6
$ws.Range("V4").Value = 0.1424305555555555
$ws.Range("W4").Value = -0.002899634147722768
$ws.Range("X4").Value = 0.07996287722380738
$ws.Range("Y4").Value = -0.08286251137153014
$ws.Range("Z4").Value = 0.9063049158926676
$ws.Range("AA4").Value = 0.007893646268240942
$ws.Range("AB4").Value = 0.06602250383131431
$ws.Range("AC4").Value = -0.05812885756307336
$ws.Range("AD4").Value = 2041.6
$ws.Range("AE4").Value = 274.35699072195
$ws.Range("AF4").Value = 2315.95699072195
$ws.Range("AG4").Value = 1085.35699072195
$ws.Range("AH4").Value = 0.2113879228152518
$ws.Range("AI4").Value = 0.2077229893121321
$ws.Range("AJ4").Value = 0.1116007352488332
$ws.Range("AK4").Value = 0.1094258014706233
$ws.Range("AL4").Value = 89.5
$ws.Range("AM4").Value = 89.5
$ws.Range("AN4").Value = 7.330700179533213
$ws.Range("AO4").Value = 1.007821229050279
$ws.Range("AP4").Value = 3.897152569917234
$ws.Range("AQ4").Value = 1.007821229050279
